$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoopFilter LPF")
$ws.Range("B2").Value = 50
